$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in this sheet are stored as plain text (shared strings), even the
# numeric-looking ones (region ids, weights). A direct `.Value = "..."` write
# lets Excel auto-coerce digit-only / decimal-looking text into a real number,
# which would also force a NumberFormat style change. To avoid both issues we
# write the literal text as a string-formula (e.g. ="2") and then flatten it
# back to a static value with a values-only Paste Special, which keeps the
# shared-string/text typing intact without touching any cell styles.

function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Write column-by-column (matches the shared-string insertion order of the
# source export): Sub.Network, then Region.Number, Region.Name,
# Sub.Network.Weight, Sub_Network.Average.of.reduced.sum for the new row 2
# content and the newly appended row 3.
Set-TextValue $ws.Range("A3") "2"
Set-TextValue $ws.Range("B3") "2"

Set-TextValue $ws.Range("C2") "18, 38, 52, 66, 72, 75, 24, 27, 33, 34, 37, 45, 65, 68, 23, 41, 25, 67, 71, 1, 5, 6, 7, 9, 14, 20, 22, 29, 30, 31"
Set-TextValue $ws.Range("C3") "59, 57"

Set-TextValue $ws.Range("D2") "Left caudalanteriorcingulate, Left posteriorcingulate, Right caudalanteriorcingulate, Right paracentral, Right posteriorcingulate, Right rostralanteriorcingulate, Left inferiortemporal, Left lateralorbitofrontal, Left parsopercularis, Left parsorbitalis, Left postcentral, Left superiortemporal, Right parahippocampal, Right parsorbitalis, Left inferiorparietal, Left rostralanteriorcingulate, Left isthmuscingulate, Right parsopercularis, Right postcentral, Left Cerebellum-Cortex, Left Pallidum, Left Hippocampus, Left Amygdala, Right Cerebellum-Cortex, Right Hippocampus, Left cuneus, Left fusiform, Left medialorbitofrontal, Left middletemporal, Left parahippocampal"
Set-TextValue $ws.Range("D3") "Right isthmuscingulate, Right inferiorparietal"

Set-TextValue $ws.Range("E2") "10.9243850276628"
Set-TextValue $ws.Range("E3") "0.0389598348666482"

Set-TextValue $ws.Range("F2") "10.9243850276628"
Set-TextValue $ws.Range("F3") "0.0389598348666482"
